$wb = $excel.ActiveWorkbook

# --- Sheet "Home win" (sheet1) ---
$ws1 = $wb.Worksheets.Item("Home win")

# Update row 2 with the new match data
$ws1.Range("A2").Value = "07-01-2025 13:00"
$ws1.Range("B2").Value = "ENGLAND"
$ws1.Range("C2").Value = "PROFESSIONAL DEVELOPMENT LEAGUE"
$ws1.Range("D2").Value = "AFC Bournemouth U21 - Bristol City U21"
$ws1.Range("E2").Value = 73.3
$ws1.Range("F2").Value = 2.25

# Remove rows 3-5 (old extra matches no longer present)
$ws1.Rows("3:5").Delete()

# --- Sheet "Draw" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Draw")

# Add a new row 3 with additional match data
$ws2.Range("A3").Value = "07-01-2025 19:00"
$ws2.Range("B3").Value = "ENGLAND"
$ws2.Range("C3").Value = "EFL TROPHY"
$ws2.Range("D3").Value = "Port Vale - Wrexham"
$ws2.Range("E3").Value = 73.3
$ws2.Range("F3").Value = 3.4
